$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ptdss1"
$ws.Cells.Item(2,3).Value = "Scarb1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 7.135072333333333
$ws.Cells.Item(2,8).Value = 21.405217
$ws.Cells.Item(2,9).Value = 0.2435182897332695
$ws.Cells.Item(2,10).Value = 0.2435182897332695
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 82.48638166666666
$ws.Cells.Item(2,14).Value = 247.459145
$ws.Cells.Item(2,15).Value = 0.7894957391680832
$ws.Cells.Item(2,16).Value = 0.7894957391680832
$ws.Cells.Item(2,17).Value = 588.5462997066072
$ws.Cells.Item(2,18).Value = 5296.916697359465
$ws.Cells.Item(2,19).Value = 0.192256652153915
$ws.Cells.Item(2,20).Value = 0.192256652153915
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ptdss1"
$ws.Cells.Item(3,3).Value = "Scarb1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 7.135072333333333
$ws.Cells.Item(3,8).Value = 21.405217
$ws.Cells.Item(3,9).Value = 0.2435182897332695
$ws.Cells.Item(3,10).Value = 0.2435182897332695
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.8713403333333334
$ws.Cells.Item(3,14).Value = 2.614021
$ws.Cells.Item(3,15).Value = 0.008339794601633706
$ws.Cells.Item(3,16).Value = 0.008339794601633706
$ws.Cells.Item(3,17).Value = 6.217076305284111
$ws.Cells.Item(3,18).Value = 55.95368674755701
$ws.Cells.Item(3,19).Value = 0.002030892518116593
$ws.Cells.Item(3,20).Value = 0.002030892518116593
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ptdss1"
$ws.Cells.Item(4,3).Value = "Scarb1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 7.135072333333333
$ws.Cells.Item(4,8).Value = 21.405217
$ws.Cells.Item(4,9).Value = 0.2435182897332695
$ws.Cells.Item(4,10).Value = 0.2435182897332695
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 17.88507033333333
$ws.Cells.Item(4,14).Value = 53.65521099999999
$ws.Cells.Item(4,15).Value = 0.1711820368112258
$ws.Cells.Item(4,16).Value = 0.1711820368112258
$ws.Cells.Item(4,17).Value = 127.6112705150874
$ws.Cells.Item(4,18).Value = 1148.501434635787
$ws.Cells.Item(4,19).Value = 0.04168595683732729
$ws.Cells.Item(4,20).Value = 0.04168595683732729
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Ptdss1"
$ws.Cells.Item(5,3).Value = "Scarb1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 7.135072333333333
$ws.Cells.Item(5,8).Value = 21.405217
$ws.Cells.Item(5,9).Value = 0.2435182897332695
$ws.Cells.Item(5,10).Value = 0.2435182897332695
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.237038999999999
$ws.Cells.Item(5,14).Value = 9.711116999999998
$ws.Cells.Item(5,15).Value = 0.03098242941905719
$ws.Cells.Item(5,16).Value = 0.03098242941905719
$ws.Cells.Item(5,17).Value = 23.09650741082099
$ws.Cells.Item(5,18).Value = 207.868566697389
$ws.Cells.Item(5,19).Value = 0.00754478822391054
$ws.Cells.Item(5,20).Value = 0.00754478822391054
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ptdss1"
$ws.Cells.Item(6,3).Value = "Scarb1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 11.272738
$ws.Cells.Item(6,8).Value = 33.818214
$ws.Cells.Item(6,9).Value = 0.3847358162785133
$ws.Cells.Item(6,10).Value = 0.3847358162785133
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 82.48638166666666
$ws.Cells.Item(6,14).Value = 247.459145
$ws.Cells.Item(6,15).Value = 0.7894957391680832
$ws.Cells.Item(6,16).Value = 0.7894957391680832
$ws.Cells.Item(6,17).Value = 929.8473690963365
$ws.Cells.Item(6,18).Value = 8368.626321867028
$ws.Cells.Item(6,19).Value = 0.3037472876572407
$ws.Cells.Item(6,20).Value = 0.3037472876572407
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ptdss1"
$ws.Cells.Item(7,3).Value = "Scarb1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 11.272738
$ws.Cells.Item(7,8).Value = 33.818214
$ws.Cells.Item(7,9).Value = 0.3847358162785133
$ws.Cells.Item(7,10).Value = 0.3847358162785133
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.8713403333333334
$ws.Cells.Item(7,14).Value = 2.614021
$ws.Cells.Item(7,15).Value = 0.008339794601633706
$ws.Cells.Item(7,16).Value = 0.008339794601633706
$ws.Cells.Item(7,17).Value = 9.822391286499332
$ws.Cells.Item(7,18).Value = 88.401521578494
$ws.Cells.Item(7,19).Value = 0.003208617683654682
$ws.Cells.Item(7,20).Value = 0.003208617683654682
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Ptdss1"
$ws.Cells.Item(8,3).Value = "Scarb1"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 11.272738
$ws.Cells.Item(8,8).Value = 33.818214
$ws.Cells.Item(8,9).Value = 0.3847358162785133
$ws.Cells.Item(8,10).Value = 0.3847358162785133
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 17.88507033333333
$ws.Cells.Item(8,14).Value = 53.65521099999999
$ws.Cells.Item(8,15).Value = 0.1711820368112258
$ws.Cells.Item(8,16).Value = 0.1711820368112258
$ws.Cells.Item(8,17).Value = 201.6137119792393
$ws.Cells.Item(8,18).Value = 1814.523407813154
$ws.Cells.Item(8,19).Value = 0.06585986066478548
$ws.Cells.Item(8,20).Value = 0.06585986066478548
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Ptdss1"
$ws.Cells.Item(9,3).Value = "Scarb1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 11.272738
$ws.Cells.Item(9,8).Value = 33.818214
$ws.Cells.Item(9,9).Value = 0.3847358162785133
$ws.Cells.Item(9,10).Value = 0.3847358162785133
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.237038999999999
$ws.Cells.Item(9,14).Value = 9.711116999999998
$ws.Cells.Item(9,15).Value = 0.03098242941905719
$ws.Cells.Item(9,16).Value = 0.03098242941905719
$ws.Cells.Item(9,17).Value = 36.49029254278199
$ws.Cells.Item(9,18).Value = 328.4126328850379
$ws.Cells.Item(9,19).Value = 0.01192005027283239
$ws.Cells.Item(9,20).Value = 0.01192005027283239
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Ptdss1"
$ws.Cells.Item(10,3).Value = "Scarb1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 6.311962000000001
$ws.Cells.Item(10,8).Value = 18.935886
$ws.Cells.Item(10,9).Value = 0.2154257335164676
$ws.Cells.Item(10,10).Value = 0.2154257335164676
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 82.48638166666666
$ws.Cells.Item(10,14).Value = 247.459145
$ws.Cells.Item(10,15).Value = 0.7894957391680832
$ws.Cells.Item(10,16).Value = 0.7894957391680832
$ws.Cells.Item(10,17).Value = 520.6509065974967
$ws.Cells.Item(10,18).Value = 4685.858159377471
$ws.Cells.Item(10,19).Value = 0.1700776987184101
$ws.Cells.Item(10,20).Value = 0.1700776987184101
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Ptdss1"
$ws.Cells.Item(11,3).Value = "Scarb1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 6.311962000000001
$ws.Cells.Item(11,8).Value = 18.935886
$ws.Cells.Item(11,9).Value = 0.2154257335164676
$ws.Cells.Item(11,10).Value = 0.2154257335164676
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.8713403333333334
$ws.Cells.Item(11,14).Value = 2.614021
$ws.Cells.Item(11,15).Value = 0.008339794601633706
$ws.Cells.Item(11,16).Value = 0.008339794601633706
$ws.Cells.Item(11,17).Value = 5.499867073067334
$ws.Cells.Item(11,18).Value = 49.49880365760601
$ws.Cells.Item(11,19).Value = 0.001796606369433618
$ws.Cells.Item(11,20).Value = 0.001796606369433617
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Ptdss1"
$ws.Cells.Item(12,3).Value = "Scarb1"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 6.311962000000001
$ws.Cells.Item(12,8).Value = 18.935886
$ws.Cells.Item(12,9).Value = 0.2154257335164676
$ws.Cells.Item(12,10).Value = 0.2154257335164676
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 17.88507033333333
$ws.Cells.Item(12,14).Value = 53.65521099999999
$ws.Cells.Item(12,15).Value = 0.1711820368112258
$ws.Cells.Item(12,16).Value = 0.1711820368112258
$ws.Cells.Item(12,17).Value = 112.8898843113273
$ws.Cells.Item(12,18).Value = 1016.008958801946
$ws.Cells.Item(12,19).Value = 0.03687701584490129
$ws.Cells.Item(12,20).Value = 0.03687701584490128
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Ptdss1"
$ws.Cells.Item(13,3).Value = "Scarb1"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 6.311962000000001
$ws.Cells.Item(13,8).Value = 18.935886
$ws.Cells.Item(13,9).Value = 0.2154257335164676
$ws.Cells.Item(13,10).Value = 0.2154257335164676
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 3.237038999999999
$ws.Cells.Item(13,14).Value = 9.711116999999998
$ws.Cells.Item(13,15).Value = 0.03098242941905719
$ws.Cells.Item(13,16).Value = 0.03098242941905719
$ws.Cells.Item(13,17).Value = 20.432067160518
$ws.Cells.Item(13,18).Value = 183.888604444662
$ws.Cells.Item(13,19).Value = 0.00667441258372258
$ws.Cells.Item(13,20).Value = 0.006674412583722579
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Ptdss1"
$ws.Cells.Item(14,3).Value = "Scarb1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 4.580172
$ws.Cells.Item(14,8).Value = 13.740516
$ws.Cells.Item(14,9).Value = 0.1563201604717497
$ws.Cells.Item(14,10).Value = 0.1563201604717497
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 82.48638166666666
$ws.Cells.Item(14,14).Value = 247.459145
$ws.Cells.Item(14,15).Value = 0.7894957391680832
$ws.Cells.Item(14,16).Value = 0.7894957391680832
$ws.Cells.Item(14,17).Value = 377.80181569098
$ws.Cells.Item(14,18).Value = 3400.21634121882
$ws.Cells.Item(14,19).Value = 0.1234141006385174
$ws.Cells.Item(14,20).Value = 0.1234141006385174
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Ptdss1"
$ws.Cells.Item(15,3).Value = "Scarb1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 4.580172
$ws.Cells.Item(15,8).Value = 13.740516
$ws.Cells.Item(15,9).Value = 0.1563201604717497
$ws.Cells.Item(15,10).Value = 0.1563201604717497
$ws.Cells.Item(15,11).Value = 2
$ws.Cells.Item(15,12).Value = 0.6666666666666666
$ws.Cells.Item(15,13).Value = 0.8713403333333334
$ws.Cells.Item(15,14).Value = 2.614021
$ws.Cells.Item(15,15).Value = 0.008339794601633706
$ws.Cells.Item(15,16).Value = 0.008339794601633706
$ws.Cells.Item(15,17).Value = 3.990888597204
$ws.Cells.Item(15,18).Value = 35.917997374836
$ws.Cells.Item(15,19).Value = 0.001303678030428813
$ws.Cells.Item(15,20).Value = 0.001303678030428813
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Ptdss1"
$ws.Cells.Item(16,3).Value = "Scarb1"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 4.580172
$ws.Cells.Item(16,8).Value = 13.740516
$ws.Cells.Item(16,9).Value = 0.1563201604717497
$ws.Cells.Item(16,10).Value = 0.1563201604717497
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 17.88507033333333
$ws.Cells.Item(16,14).Value = 53.65521099999999
$ws.Cells.Item(16,15).Value = 0.1711820368112258
$ws.Cells.Item(16,16).Value = 0.1711820368112258
$ws.Cells.Item(16,17).Value = 81.91669835876399
$ws.Cells.Item(16,18).Value = 737.2502852288759
$ws.Cells.Item(16,19).Value = 0.02675920346421179
$ws.Cells.Item(16,20).Value = 0.02675920346421179
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Ptdss1"
$ws.Cells.Item(17,3).Value = "Scarb1"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 4.580172
$ws.Cells.Item(17,8).Value = 13.740516
$ws.Cells.Item(17,9).Value = 0.1563201604717497
$ws.Cells.Item(17,10).Value = 0.1563201604717497
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 3.237038999999999
$ws.Cells.Item(17,14).Value = 9.711116999999998
$ws.Cells.Item(17,15).Value = 0.03098242941905719
$ws.Cells.Item(17,16).Value = 0.03098242941905719
$ws.Cells.Item(17,17).Value = 14.826195390708
$ws.Cells.Item(17,18).Value = 133.435758516372
$ws.Cells.Item(17,19).Value = 0.004843178338591678
$ws.Cells.Item(17,20).Value = 0.004843178338591677
